$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = -0.060437723747241136
$ws.Range("B1").Value = 0.06043772172119339

$ws.Range("A2").Value = 0.013116550131572306
$ws.Range("B2").Value = -0.013116552290001597

$ws.Range("A3").Value = 0.053351735473098724
$ws.Range("B3").Value = -0.05335173757708405

$ws.Range("A4").Value = -0.077494404745110906
$ws.Range("B4").Value = 0.07749440273053454

$ws.Range("A5").Value = 0.0092644102004602255
$ws.Range("B5").Value = -0.0092644122988238473
